# Restore C10 on the "Rules" sheet (the "From" value for rule R30)
# from 18 to 1, as captured by the commit's revision restore.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
